$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C header (13-01-2023), matching B1 style (bold/border/center)
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "13-01-2023"

$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 555521.72
$ws.Range("C2").Value = 780789.63
$ws.Range("A3").Value = "1822 Raices Valores Negociables"
$ws.Range("B3").Value = 1965552.25
$ws.Range("C3").Value = 1964743.29
$ws.Range("A4").Value = "Adcap IOL Acciones Argentina"
$ws.Range("B4").Value = 180458.44
$ws.Range("C4").Value = 191274.98
$ws.Range("A5").Value = "Allaria Acciones"
$ws.Range("B5").Value = 409772.54
$ws.Range("C5").Value = 432395.27
$ws.Range("A6").Value = "Alpha Acciones"
$ws.Range("B6").Value = 641526.89
$ws.Range("C6").Value = 641456.41
$ws.Range("A7").Value = "Alpha Mega"
$ws.Range("B7").Value = 2121018
$ws.Range("C7").Value = 2121140.19
$ws.Range("A8").Value = "Alpha Mercosur"
$ws.Range("B8").Value = 60203.17
$ws.Range("C8").Value = 60345.62
$ws.Range("A9").Value = "Alpha Recursos Naturales"
$ws.Range("B9").Value = 280549.07
$ws.Range("C9").Value = 315455.6
$ws.Range("A10").Value = "Alpha planeam equil"
$ws.Range("B10").Value = 22182.47
$ws.Range("C10").Value = 14634.87
$ws.Range("A11").Value = "Alpha renta balan global"
$ws.Range("B11").Value = 483598.22
$ws.Range("C11").Value = 482059.93
$ws.Range("A12").Value = "Argenfunds"
$ws.Range("B12").Value = 59321.84
$ws.Range("C12").Value = 59322.28
$ws.Range("A13").Value = "Arpenta ex Mercosur"
$ws.Range("B13").Value = 13016.39
$ws.Range("C13").Value = 13039.74
$ws.Range("A14").Value = "Balanz"
$ws.Range("B14").Value = 1609557.2
$ws.Range("C14").Value = 1675673.22
$ws.Range("A15").Value = "Bull Market"
$ws.Range("B15").Value = 190921.05
$ws.Range("C15").Value = 190918.76
$ws.Range("A16").Value = "CMA acciones"
$ws.Range("B16").Value = 434655.43
$ws.Range("C16").Value = 434722.52
$ws.Range("A17").Value = "Compass Crecimiento"
$ws.Range("B17").Value = 5520011.45
$ws.Range("C17").Value = 5478700.07
$ws.Range("A18").Value = "Compass Crecimiento II"
$ws.Range("B18").Value = 40199.72
$ws.Range("C18").Value = 40202.68
$ws.Range("A19").Value = "Consultatio Acciones Argentina"
$ws.Range("B19").Value = 3718020.23
$ws.Range("C19").Value = 3622166.41
$ws.Range("A20").Value = "Consultatio Renta Variable"
$ws.Range("B20").Value = 1128360.79
$ws.Range("C20").Value = 1127816.33
$ws.Range("A21").Value = "Delta Acciones"
$ws.Range("B21").Value = 500230.26
$ws.Range("C21").Value = 500228.56
$ws.Range("A22").Value = "Delta Internacional"
$ws.Range("B22").Value = 19993.64
$ws.Range("C22").Value = 20007.35
$ws.Range("A23").Value = "Delta Latinoamerica"
$ws.Range("B23").Value = 45407.82
$ws.Range("C23").Value = 45414.76
$ws.Range("A24").Value = "Delta Select"
$ws.Range("B24").Value = 3101616.47
$ws.Range("C24").Value = 3049799.52
$ws.Range("A25").Value = "Delta gestion V"
$ws.Range("B25").Value = 297494.04
$ws.Range("C25").Value = 298120.48
$ws.Range("A26").Value = "FBA Acciones Argentinas"
$ws.Range("B26").Value = 1409067.53
$ws.Range("C26").Value = 1438504.19
$ws.Range("A27").Value = "FBA Calificado"
$ws.Range("B27").Value = 1386031.58
$ws.Range("C27").Value = 1411582.04
$ws.Range("A28").Value = "Fima Acciones"
$ws.Range("B28").Value = 1860247.18
$ws.Range("C28").Value = 2201916.87
$ws.Range("A29").Value = "Fima PB Acciones"
$ws.Range("B29").Value = 1084903.44
$ws.Range("C29").Value = 1317442.24
$ws.Range("A30").Value = "Gainvest Renta Variable"
$ws.Range("B30").Value = 304076
$ws.Range("C30").Value = 303875.61
$ws.Range("A31").Value = "Galileo Acciones"
$ws.Range("B31").Value = 5640649.97
$ws.Range("C31").Value = 5641847.99
$ws.Range("A32").Value = "Goal Acciones Argentinas"
$ws.Range("B32").Value = 222077.87
$ws.Range("C32").Value = 237028.35
$ws.Range("A33").Value = "Goal acciones plus"
$ws.Range("B33").Value = 50979.46
$ws.Range("C33").Value = 50988.11
$ws.Range("A34").Value = "HF Acciones Argentinas"
$ws.Range("B34").Value = 597986.57
$ws.Range("C34").Value = 587771.51
$ws.Range("A35").Value = "HF Acciones Lideres"
$ws.Range("B35").Value = 1384726.83
$ws.Range("C35").Value = 1395133.34
$ws.Range("A36").Value = "IAM Renta Variable"
$ws.Range("B36").Value = 191837.98
$ws.Range("C36").Value = 200806.93
$ws.Range("A37").Value = "IEB Value"
$ws.Range("B37").Value = 31306.2
$ws.Range("C37").Value = 31316.24
$ws.Range("A38").Value = "Lombardi"
$ws.Range("B38").Value = 261735.88
$ws.Range("C38").Value = 277947.33
$ws.Range("A39").Value = "MAF"
$ws.Range("B39").Value = 212849.36
$ws.Range("C39").Value = 212751.82
$ws.Range("A40").Value = "Megainver"
$ws.Range("B40").Value = 159198.75
$ws.Range("C40").Value = 159108.49
$ws.Range("A41").Value = "Pellegrini Acciones"
$ws.Range("B41").Value = 569589.4
$ws.Range("C41").Value = 569699.02
$ws.Range("A42").Value = "Pionero Acciones"
$ws.Range("B42").Value = 1153255.7
$ws.Range("C42").Value = 1153463.4
$ws.Range("A43").Value = "Premier Renta Variable"
$ws.Range("B43").Value = 341306.89
$ws.Range("C43").Value = 341388.75
$ws.Range("A44").Value = "Quinquela Acciones"
$ws.Range("B44").Value = 507296.33
$ws.Range("C44").Value = 507295.27
$ws.Range("A45").Value = "Rofex 20 Renta Variable"
$ws.Range("B45").Value = 367543.33
$ws.Range("C45").Value = 367494.79
$ws.Range("A46").Value = "SBS Acciones Argentina"
$ws.Range("B46").Value = 2313501.55
$ws.Range("C46").Value = 2313158.64
$ws.Range("A47").Value = "Schroeder RV"
$ws.Range("B47").Value = 3982233.18
$ws.Range("C47").Value = 3981960.68
$ws.Range("A48").Value = "Supefondo RV"
$ws.Range("B48").Value = 5890140.48
$ws.Range("C48").Value = 7475061.99
$ws.Range("A49").Value = "Superfondo "
$ws.Range("B49").Value = 42883.67
$ws.Range("C49").Value = 41495.97
$ws.Range("A50").Value = "Toronto Trust Multimercado"
$ws.Range("B50").Value = 190817.94
$ws.Range("C50").Value = 190746.08
$ws.Range("A51").Value = "Toronto trust Argy"
$ws.Range("B51").Value = 88525.25
$ws.Range("C51").Value = 88448.78
$ws.Range("A52").Value = "avg"
$ws.Range("B52").Value = 1072879.15
$ws.Range("C52").Value = 1121173.26
$ws.Range("A53").Value = "total"
$ws.Range("B53").Value = 53643957.42
$ws.Range("C53").Value = 56058662.9
